$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet 1 "User Rights": append three new user-right rows (96-98)
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

$newRights = @(
    @{
        Row = 96
        Name = "LINE_LISTING_CONFIGURE"
        Values = @("Yes","Yes","Yes","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No")
    },
    @{
        Row = 97
        Name = "LINE_LISTING_CONFIGURE_NATION"
        Values = @("Yes","Yes","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No")
    },
    @{
        Row = 98
        Name = "AGGREGATE_REPORT_VIEW"
        Values = @("Yes","Yes","Yes","Yes","Yes","No","Yes","No","Yes","No","No","No","No","Yes","Yes","Yes","Yes","Yes","Yes","Yes")
    }
)

$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

foreach ($right in $newRights) {
    $r = $right.Row

    # Column A - bold "User Right" style (copy formatting from row 95, col A)
    $ws.Cells.Item($r, 1).Value = $right.Name
    $ws.Range("A95").Copy()
    $ws.Cells.Item($r, 1).PasteSpecial($xlPasteFormats)

    # Column B - plain duplicate of the right name (copy formatting from row 95, col B)
    $ws.Cells.Item($r, 2).Value = $right.Name
    $ws.Range("B95").Copy()
    $ws.Cells.Item($r, 2).PasteSpecial($xlPasteFormats)

    # Columns C..V - Yes/No flags per role
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $val = $right.Values[$i]
        $cell = $ws.Range($col + $r)
        $cell.Value = $val
        if ($val -eq "Yes") {
            $ws.Range("C95").Copy()
        } else {
            $ws.Range("D95").Copy()
        }
        $cell.PasteSpecial($xlPasteFormats)
    }
}

# ----------------------------------------------------------------------
# Sheet 2 "About": bump the SORMAS version string
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").Value = "1.31.0-SNAPSHOT"
